$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format temporarily so numeric-looking strings
# (e.g. "88.01", "0.539") are stored as text, matching the source data,
# then restore the original "Normal" style so no stray formatting remains.
$dCol = $ws.Range("D2:D50")
$dCol.NumberFormat = "@"

$ws.Range("D2").Value = "69.701.77"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").Value = "3.749.23"
$ws.Range("E3").Value = "  +7.19%  "

$ws.Range("D5").Value = "612.47"
$ws.Range("E5").Value = "  +4.29%  "

$ws.Range("D6").Value = "177.78"
$ws.Range("E6").Value = "  -2.55%  "

$ws.Range("D7").Value = "3.746.04"
$ws.Range("E7").Value = "  +7.39%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  +1.58%  "

$ws.Range("E10").Value = "  +6.10%  "

$ws.Range("D11").Value = "6.33"
$ws.Range("E11").Value = "  -3.55%  "

$ws.Range("D12").Value = "0.496"
$ws.Range("E12").Value = "  +1.76%  "

$ws.Range("D13").Value = "40.78"
$ws.Range("E13").Value = "  +6.61%  "

$ws.Range("E14").Value = "  +2.02%  "

$ws.Range("D15").Value = "4.371.17"
$ws.Range("E15").Value = "  +7.07%  "

$ws.Range("D16").Value = "3.746.54"
$ws.Range("E16").Value = "  +6.97%  "

$ws.Range("D17").Value = "69.818.82"
$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("E18").Value = "  +0.56%  "

$ws.Range("D19").Value = "7.59"
$ws.Range("E19").Value = "  +2.68%  "

$ws.Range("D20").Value = "514.96"
$ws.Range("E20").Value = "  +2.27%  "

$ws.Range("D21").Value = "16.67"
$ws.Range("E21").Value = "  -0.67%  "

$ws.Range("D22").Value = "9.53"
$ws.Range("E22").Value = "  +6.89%  "

$ws.Range("D23").Value = "0.726"
$ws.Range("E23").Value = "  -0.48%  "

$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D24").Value = "2.52"
$ws.Range("E24").Value = "  +6.34%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "88.01"
$ws.Range("E25").Value = "  +1.84%  "

$ws.Range("D26").Value = "13.35"
$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("D27").Value = "11.06"
$ws.Range("E27").Value = "  +3.62%  "

$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("E29").Value = "  +17.67%  "

$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "2.84"
$ws.Range("E31").Value = "  +4.83%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "7.84"
$ws.Range("E32").Value = "  -2.77%  "

$ws.Range("D33").Value = "31.34"
$ws.Range("E33").Value = "  +2.66%  "

$ws.Range("E34").Value = "  -0.64%  "

$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("E36").Value = "  +2.62%  "

$ws.Range("D37").Value = "1.04"
$ws.Range("E37").Value = "  +2.71%  "

$ws.Range("E38").Value = "  +2.87%  "

$ws.Range("E39").Value = "  +3.24%  "

$ws.Range("E40").Value = "  +4.93%  "

$ws.Range("D41").Value = "51.17"
$ws.Range("E41").Value = "  +1.78%  "

$ws.Range("D42").Value = "44.37"
$ws.Range("E42").Value = "  -4.80%  "

$ws.Range("D43").Value = "8.81"
$ws.Range("E43").Value = "  +1.68%  "

$ws.Range("D44").Value = "422.48"
$ws.Range("E44").Value = "  +4.21%  "

$ws.Range("D45").Value = "3.078.51"
$ws.Range("E45").Value = "  +2.92%  "

$ws.Range("E46").Value = "  -1.82%  "

$ws.Range("D47").Value = "0.0364"
$ws.Range("E47").Value = "  +0.47%  "

$ws.Range("D48").Value = "27.77"
$ws.Range("E48").Value = "  +0.47%  "

$ws.Range("D49").Value = "2.52"
$ws.Range("E49").Value = "  +3.97%  "

$ws.Range("D50").Value = "136.13"
$ws.Range("E50").Value = "  +1.00%  "

# Restore default styling on column D (undo temporary text format)
$dCol.Style = "Normal"
